# Update "国内生产总值环比增长速度" sheet:
# the B/C quarter rows (and their values) for each year were reordered,
# and several growth-rate figures were revised.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 2011
$ws.Range("B2").Value = 2.3
$ws.Range("A3").Value = "2011年C"
$ws.Range("B3").Value = 1.8
$ws.Range("A4").Value = "2011年B"
$ws.Range("B4").Value = 2.3

# 2012
$ws.Range("B6").Value = 1.7
$ws.Range("A7").Value = "2012年C"
$ws.Range("B7").Value = 1.9
$ws.Range("A8").Value = "2012年B"
$ws.Range("B8").Value = 2.2

# 2013
$ws.Range("B10").Value = 1.7
$ws.Range("A11").Value = "2013年C"
$ws.Range("B11").Value = 2.2
$ws.Range("A12").Value = "2013年B"
$ws.Range("B12").Value = 1.8

# 2014
$ws.Range("B14").Value = 1.7
$ws.Range("A15").Value = "2014年C"
$ws.Range("A16").Value = "2014年B"
$ws.Range("B16").Value = 1.9
$ws.Range("B17").Value = 1.7

# 2015
$ws.Range("B18").Value = 1.6
$ws.Range("A19").Value = "2015年C"
$ws.Range("B19").Value = 1.7
$ws.Range("A20").Value = "2015年B"
$ws.Range("B20").Value = 2

# 2016
$ws.Range("B21").Value = 1.6
$ws.Range("B22").Value = 1.6
$ws.Range("A23").Value = "2016年C"
$ws.Range("B23").Value = 1.7
$ws.Range("A24").Value = "2016年B"
$ws.Range("B24").Value = 1.8
$ws.Range("B25").Value = 1.5

# 2017
$ws.Range("A27").Value = "2017年C"
$ws.Range("B27").Value = 1.6
$ws.Range("A28").Value = "2017年B"
$ws.Range("B28").Value = 1.8
$ws.Range("B29").Value = 1.5

# 2018
$ws.Range("A31").Value = "2018年C"
$ws.Range("B31").Value = 1.5
$ws.Range("A32").Value = "2018年B"
$ws.Range("B32").Value = 1.7
$ws.Range("B33").Value = 1.3

# 2019
$ws.Range("A35").Value = "2019年C"
$ws.Range("B35").Value = 1.5
$ws.Range("A36").Value = "2019年B"
$ws.Range("B36").Value = 1.5
$ws.Range("B37").Value = 1.2

# 2020
$ws.Range("B38").Value = -10.4
$ws.Range("A39").Value = "2020年C"
$ws.Range("B39").Value = 3.5
$ws.Range("A40").Value = "2020年B"
$ws.Range("B40").Value = 11.5
$ws.Range("B41").Value = 2.5

# 2021
$ws.Range("B42").Value = 0.5
$ws.Range("A43").Value = "2021年C"
$ws.Range("B43").Value = 0.7
$ws.Range("A44").Value = "2021年B"
$ws.Range("B44").Value = 1.3
